# Apply the "new version with timestamp" update to the daily sales report.
#
# Summary of the change (per the author's diff):
#  - H11 ("مسك الرمان" item's time range) changes from "19:0" to "17:0"
#  - P11 (that item's sale price) changes from "30.0000" to "90.0000"
#    (still stored as text, matching the original cell formatting)
#  - Q11 (that item's transaction count) changes to a new value "3:0"
#  - P12 (the sale-price column total) is recalculated to reflect the
#    updated P11 value: 294.22000000000003 -> 354.22000000000003

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H11: "19:0" -> "17:0" (plain text cell, format already Text) ---
$h11 = $ws.Range("H11")
$h11.Value = "17:0"

# --- P11: "30.0000" -> "90.0000" ---
# This cell's number format is a numeric one (0.00) even though the
# stored value is textual. Assigning a numeric-looking string directly
# would make Excel coerce it into a real number, so we temporarily force
# a text format, assign the text, then restore the original number
# format (the underlying value stays text, matching the source file).
$p11 = $ws.Range("P11")
$p11OriginalFormat = $p11.NumberFormat()
$p11.NumberFormat = "@"
$p11.Value = "90.0000"
$p11.NumberFormat = $p11OriginalFormat

# --- Q11: new transaction-count value "3:0" ---
$q11 = $ws.Range("Q11")
$q11.Value = "3:0"

# --- P12: recalculate the sale-price total ---
$p12 = $ws.Range("P12")
$p12.Value = 354.22000000000003
